$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update RF column (I) for rows 22-58 with the new recalculated factor
$ws.Range("I22:I58").Value = 13.31555555555556
